$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: suite finished automating -> status changes from "Testing" to "Automated",
# and case counts are updated (Automated: 4 -> 7, Total: 6 -> 7)
$ws.Range("D3").Value = "Automated"
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 7

# Update the active selection to E3
$ws.Activate()
$ws.Range("E3").Select()
